# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-27
$kValues = @{
    2  = 7
    3  = 4
    4  = 6
    5  = 6
    6  = 9
    7  = 2
    8  = 5
    9  = 2
    10 = 6
    11 = 3
    12 = 6
    13 = 7
    14 = 5
    15 = 8
    16 = 5
    17 = 11
    18 = 8
    19 = 3
    20 = 11
    21 = 7
    22 = 3
    23 = 6
    24 = 10
    25 = 6
    26 = 7
    27 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
